# Auto-generated Excel COM-interop edit script
# Refreshes the cryptocurrency price/volume table (scheduled GitHub Actions scrape).
#
# Price values in column D are plain text in the source sheet (some look like
# thousands-grouped numbers, e.g. '34.547.83', which are NOT valid numbers; others
# look like plain decimals, e.g. '230.06', which Excel's normal cell-input parsing
# WOULD auto-convert to a number). Set-TextValue below forces text the same way a
# user would in the Excel UI (leading apostrophe), then restores the cell's original
# style so no visible/number formatting changes are introduced.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $c = $ws.Range($addr)
    $origStyle = $c.Style
    $c.Value = "'" + $val
    $c.Style = $origStyle
}

$ws.Range("D2").Value = '34.547.83'
$ws.Range("E2").Value = '  -2.72%  '

$ws.Range("D3").Value = '1.805.51'
$ws.Range("E3").Value = '  -1.90%  '

$ws.Range("E4").Value = '  +0.55%  '

Set-TextValue "D5" '230.06'
$ws.Range("E5").Value = '  -0.65%  '

Set-TextValue "D6" '0.608'
$ws.Range("E6").Value = '  -0.41%  '

$ws.Range("E7").Value = '  +0.59%  '

Set-TextValue "D8" '39.01'
$ws.Range("E8").Value = '  -11.24%  '

Set-TextValue "D9" '0.321'
$ws.Range("E9").Value = '  +3.08%  '

Set-TextValue "D10" '0.0678'
$ws.Range("E10").Value = '  -3.57%  '

$ws.Range("E11").Value = '  -1.82%  '

$ws.Range("D12").Value = '2.068.47'
$ws.Range("E12").Value = '  -1.82%  '

# Rows 13-15: coin rows rotated (WrappedEther/Chainlink/Polygon -> Chainlink/Polygon/WrappedEther)
$ws.Range("B13").Value = 'Chainlink'
$ws.Range("C13").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextValue "D13" '11.13'
$ws.Range("E13").Value = '  -1.15%  '

$ws.Range("B14").Value = 'Polygon'
$ws.Range("C14").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
Set-TextValue "D14" '0.659'
$ws.Range("E14").Value = '  -2.49%  '

$ws.Range("B15").Value = 'WrappedEther'
$ws.Range("C15").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D15").Value = '1.775.91'
$ws.Range("E15").Value = '  -3.44%  '

Set-TextValue "D16" '4.55'
$ws.Range("E16").Value = '  -4.31%  '

$ws.Range("D17").Value = '34.592.28'
$ws.Range("E17").Value = '  -2.50%  '

Set-TextValue "D18" '68.95'
$ws.Range("E18").Value = '  -2.00%  '

Set-TextValue "D19" '242.68'
$ws.Range("E19").Value = '  -0.77%  '

$ws.Range("D20").Value = '0.0₃0778'
$ws.Range("E20").Value = '  -3.02%  '

Set-TextValue "D21" '11.77'
$ws.Range("E21").Value = '  -2.35%  '

$ws.Range("E22").Value = '  -1.22%  '

$ws.Range("E23").Value = '  +0.50%  '

Set-TextValue "D24" '2.23'
$ws.Range("E24").Value = '  +0.15%  '

Set-TextValue "D25" '172.27'
$ws.Range("E25").Value = '  +0.42%  '

Set-TextValue "D26" '7.72'
$ws.Range("E26").Value = '  -3.21%  '

Set-TextValue "D27" '17.14'
$ws.Range("E27").Value = '  -3.90%  '

$ws.Range("E28").Value = '  -0.68%  '

$ws.Range("E29").Value = '  -5.86%  '

$ws.Range("E30").Value = '  +0.65%  '

Set-TextValue "D31" '4.03'
$ws.Range("E31").Value = '  +2.21%  '

Set-TextValue "D32" '0.0540'
$ws.Range("E32").Value = '  -2.51%  '

Set-TextValue "D33" '3.88'
$ws.Range("E33").Value = '  -5.02%  '

$ws.Range("E34").Value = '  +10.23%  '

Set-TextValue "D35" '1.78'
$ws.Range("E35").Value = '  -4.07%  '

Set-TextValue "D36" '0.690'
$ws.Range("E36").Value = '  -0.36%  '

Set-TextValue "D37" '90.89'
$ws.Range("E37").Value = '  -5.16%  '

$ws.Range("E38").Value = '  +4.43%  '

$ws.Range("D39").Value = '1.319.95'
$ws.Range("E39").Value = '  -2.30%  '

$ws.Range("E40").Value = '  -2.45%  '

$ws.Range("E41").Value = '  +0.04%  '

Set-TextValue "D42" '0.951'
$ws.Range("E42").Value = '  -6.26%  '

Set-TextValue "D43" '14.24'
$ws.Range("E43").Value = '  -8.92%  '

$ws.Range("E44").Value = '  -3.57%  '

Set-TextValue "D45" '2.19'
$ws.Range("E45").Value = '  -10.87%  '

Set-TextValue "D46" '6.18'
$ws.Range("E46").Value = '  -1.93%  '

Set-TextValue "D47" '0.0513'
$ws.Range("E47").Value = '  -1.14%  '

$ws.Range("D48").Value = '1.990.36'
$ws.Range("E48").Value = '  -0.94%  '

$ws.Range("E49").Value = '  +0.57%  '

Set-TextValue "D50" '0.0662'
$ws.Range("E50").Value = '  +4.56%  '

Set-TextValue "D51" '97.53'
$ws.Range("E51").Value = '  -5.88%  '
